# Added learnings on 17/09/2017
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10: change "Method Overriding" -> "Generics" and add description in B10
$ws.Range("A10").Value = "Generics "
$ws.Range("B10").Value = "It was introduced on 1.5 to provide compile-time type checking and removing risk of ClassCastException during run time."

# Column A labels for new rows
$ws.Range("A11").Value = "Autoboxing"
$ws.Range("A12").Value = "Unboxing"
$ws.Range("A13").Value = "Generic Type"

# Column B descriptions for new rows
$ws.Range("B13").Value = "A class or interface that is parameterized over types."
$ws.Range("B11").Value = "Convert primitive data types to corresponding Wrapper classes"
$ws.Range("B12").Value = "Convert  Wrapper classes to corresponding primitive data types "

# Column C examples for new rows (long, wrapped text)
$ws.Range("C12").Value = "Pass as a parameter to a method that expects a value of the corresponding primitive type.
Assign to a variable of the corresponding primitive type.
Inside main method
Integer wrapped = new Integer(100);
call(wrapped);
private static void call(int primitive) {
}"

$ws.Range("C11").Value = "Pass as a parameter to a method that expects an object of the corresponding wrapper class. For example a method with Integer argument can be called by passing int, java compiler will do the conversion of int to Integer.
Assign to a variable of the corresponding wrapper class. For example, assigning a Long object to long variable.
Inside main method
int primitive = 100;
call(primitive);
private static void call(Integer wrapped) {
int primitive  = wrapped;
wrapper = new Integer(primitive);
wrapped = primitive;
}"

$ws.Range("C11:C12").WrapText = $true
$ws.Range("C11:C12").VerticalAlignment = -4160

# Row heights for the new wrapped text rows
$ws.Rows.Item(11).RowHeight = 240
$ws.Rows.Item(12).RowHeight = 150

# Update the selection to match where the author ended up editing
$ws.Range("C11").Select()
